# Insert a new data row at row 7 (pushes existing rows 7-21 down to 8-22)
# and populate it with the new weekly price observation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(7).Insert()

$ws.Cells.Item(7, 1).Value = 11
$ws.Cells.Item(7, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(7, 3).Value = "Bíobío"
$ws.Cells.Item(7, 4).Value = 45037
$ws.Cells.Item(7, 5).Value = 8
$ws.Cells.Item(7, 6).Value = "Fruta"
$ws.Cells.Item(7, 7).Value = 100104
$ws.Cells.Item(7, 8).Value = "Frutos de pepita"
$ws.Cells.Item(7, 9).Value = 100104003
$ws.Cells.Item(7, 10).Value = "Membrillo"
$ws.Cells.Item(7, 11).Value = "Champion"
$ws.Cells.Item(7, 12).Value = "Primera"
$ws.Cells.Item(7, 13).Value = 250
$ws.Cells.Item(7, 14).Value = 9000
$ws.Cells.Item(7, 15).Value = 9500
$ws.Cells.Item(7, 16).Value = 9200
$ws.Cells.Item(7, 17).Value = "`$/caja 18 kilos granel"
$ws.Cells.Item(7, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(7, 19).Value = 511
$ws.Cells.Item(7, 20).Value = 18
